# Edit slide 1's "Presented by" credits textbox:
#  - reposition the textbox
#  - collapse "Presented by" / "Deekshita Athreya (22BCE062)" / "Vashita Darji(22BCE056)" /
#    "Ansh Bhavsar (22BCE019)" / "Sezan Agvan (22BCE012)" / (blank bullet) into a single
#    paragraph reading "Presented by: Deekshita Athreya (22BCE062)"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(60)   # "Google Shape;204;p16"

# --- reposition (EMU -> points; 12700 EMU per point) ---
$sh.Left = 197589 / 12700
$sh.Top  = 4241089 / 12700

# --- rewrite the text body as a single paragraph ---
$tr = $sh.TextFrame.TextRange
$tr.Text = "Presented by: Deekshita Athreya (22BCE062)"

# "Presented " / "by: " -> bold + underlined (as the old "Presented by" run was)
$r1 = $tr.Characters(1, 10)
$r1.Font.Bold = $true
$r1.Font.Underline = $true
$r1.Font.Name = "Roboto"

$r2 = $tr.Characters(11, 4)
$r2.Font.Bold = $true
$r2.Font.Underline = $true
$r2.Font.Name = "Roboto"

# "Deekshita Athreya (22BCE062)" keeps the plain (non-bold, non-underlined) styling
# that it already had as its own paragraph before the merge. Touch each former run's
# span individually (even where formatting matches its neighbour) to keep the same
# run boundaries as the source paragraph.
$r3 = $tr.Characters(15, 9)     # "Deekshita"
$r3.Font.Bold = $false
$r3.Font.Underline = $false
$r3.Font.Name = "Roboto"

$r4 = $tr.Characters(24, 1)     # " "
$r4.Font.Bold = $false
$r4.Font.Underline = $false
$r4.Font.Name = "Roboto"

$r5 = $tr.Characters(25, 7)     # "Athreya"
$r5.Font.Bold = $false
$r5.Font.Underline = $false
$r5.Font.Name = "Roboto"

$r6 = $tr.Characters(32, 1)     # " "
$r6.Font.Bold = $false
$r6.Font.Underline = $false
$r6.Font.Name = "Roboto"

$r7 = $tr.Characters(33, 9)     # "(22BCE062"
$r7.Font.Bold = $false
$r7.Font.Underline = $false
$r7.Font.Name = "Roboto"

$r8 = $tr.Characters(42, 1)     # ")"
$r8.Font.Bold = $false
$r8.Font.Underline = $false
$r8.Font.Name = "Roboto"
